# Daily attendance processing - 2025-10-28 19:42:07
#
# The "Recorded By" column (column G) lists who recorded/edited each
# attendance entry as a comma-separated string, e.g. "System, user@example.com".
# Whenever "System" is the first (leftmost) entry in that list, it should be
# moved so it no longer leads the list - i.e. the last recorder in the list
# is rotated to the front, pushing "System" (and everything else) back by
# one position. Lists that don't start with "System", or that consist of a
# single entry, are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $val = $cell.Value()

    if ($val -eq $null) {
        continue
    }

    $parts = $val -split ', '

    if ($parts.Length -gt 1 -and $parts[0] -eq 'System') {
        $rotated = @($parts[-1]) + $parts[0..($parts.Length - 2)]
        $newVal = $rotated -join ', '
        if ($newVal -ne $val) {
            $cell.Value = $newVal
        }
    }
}
